{"js": "// This document has a title paragraph (the date heading) followed by a single\n// 20-row x 5-column table. Only data rows 1, 5, 9, 13, 17 (1-based) contain text;\n// the rows between them are intentionally blank spacer rows.\n\n// 1) Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.getRange().insertText(\"2026-01-25 Sunday\", Word.InsertLocation.replace); // '2026-01-24 Saturday' -> '2026-01-25 Sunday'\n\n// 2) Update the division-problem answers in the table, cell by cell, so that each\n//    run keeps its existing font/size formatting (we replace just the paragraph's\n//    text range instead of clearing/recreating the cell body).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// [rowIndex (0-based), columnIndex (0-based), newText]\nconst cellUpdates = [\n  [0, 0, \"42\u00f79=4, 6\"], // '31\u00f76=5, 1' -> '42\u00f79=4, 6'\n  [0, 1, \"61\u00f73=20, 1\"], // '10\u00f74=2, 2' -> '61\u00f73=20, 1'\n  [0, 2, \"86\u00f78=10, 6\"], // '45\u00f76=7, 3' -> '86\u00f78=10, 6'\n  [0, 3, \"73\u00f72=36, 1\"], // '64\u00f79=7, 1' -> '73\u00f72=36, 1'\n  [0, 4, \"66\u00f78=8, 2\"], // '16\u00f73=5, 1' -> '66\u00f78=8, 2'\n  [4, 0, \"65\u00f78=8, 1\"], // '43\u00f72=21, 1' -> '65\u00f78=8, 1'\n  [4, 1, \"40\u00f78=5, 0\"], // '59\u00f74=14, 3' -> '40\u00f78=5, 0'\n  [4, 2, \"19\u00f79=2, 1\"], // '32\u00f75=6, 2' -> '19\u00f79=2, 1'\n  [4, 3, \"22\u00f76=3, 4\"], // '55\u00f73=18, 1' -> '22\u00f76=3, 4'\n  [4, 4, \"63\u00f76=10, 3\"], // '53\u00f74=13, 1' -> '63\u00f76=10, 3'\n  [8, 0, \"18\u00f79=2, 0\"], // '68\u00f74=17, 0' -> '18\u00f79=2, 0'\n  [8, 1, \"29\u00f79=3, 2\"], // '86\u00f72=43, 0' -> '29\u00f79=3, 2'\n  [8, 2, \"32\u00f78=4, 0\"], // '85\u00f76=14, 1' -> '32\u00f78=4, 0'\n  [8, 3, \"68\u00f73=22, 2\"], // '99\u00f78=12, 3' -> '68\u00f73=22, 2'\n  [8, 4, \"83\u00f72=41, 1\"], // '53\u00f78=6, 5' -> '83\u00f72=41, 1'\n  [12, 0, \"99\u00f78=12, 3\"], // '60\u00f73=20, 0' -> '99\u00f78=12, 3'\n  [12, 1, \"98\u00f79=10, 8\"], // '25\u00f76=4, 1' -> '98\u00f79=10, 8'\n  [12, 2, \"11\u00f76=1, 5\"], // '10\u00f75=2, 0' -> '11\u00f76=1, 5'\n  [12, 3, \"44\u00f74=11, 0\"], // '57\u00f75=11, 2' -> '44\u00f74=11, 0'\n  [12, 4, \"15\u00f76=2, 3\"], // '85\u00f76=14, 1' -> '15\u00f76=2, 3'\n  [16, 0, \"27\u00f79=3, 0\"], // '20\u00f78=2, 4' -> '27\u00f79=3, 0'\n  [16, 1, \"53\u00f73=17, 2\"], // '13\u00f76=2, 1' -> '53\u00f73=17, 2'\n  [16, 2, \"50\u00f72=25, 0\"], // '47\u00f75=9, 2' -> '50\u00f72=25, 0'\n  [16, 3, \"81\u00f73=27, 0\"], // '24\u00f78=3, 0' -> '81\u00f73=27, 0'\n  [16, 4, \"71\u00f78=8, 7\"], // '28\u00f72=14, 0' -> '71\u00f78=8, 7'\n];\n\nfor (const [rowIndex, colIndex, newText] of cellUpdates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const cellPara = cell.body.paragraphs.getFirst();\n  cellPara.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# This document has a title paragraph (the date heading) followed by a single\n# 20-row x 5-column table. Only data rows 1, 5, 9, 13, 17 contain text; the rows\n# between them are intentionally blank spacer rows.\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-25 Sunday\" # '2026-01-24 Saturday' -> '2026-01-25 Sunday'\n\n# 2) Update the division-problem answers in the table, cell by cell, so each run\n#    keeps its existing font/size formatting (Range.Text only replaces the text,\n#    it does not touch the surrounding run/paragraph properties).\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"42\u00f79=4, 6\" # '31\u00f76=5, 1' -> '42\u00f79=4, 6'\n$t.Cell(1, 2).Range.Text = \"61\u00f73=20, 1\" # '10\u00f74=2, 2' -> '61\u00f73=20, 1'\n$t.Cell(1, 3).Range.Text = \"86\u00f78=10, 6\" # '45\u00f76=7, 3' -> '86\u00f78=10, 6'\n$t.Cell(1, 4).Range.Text = \"73\u00f72=36, 1\" # '64\u00f79=7, 1' -> '73\u00f72=36, 1'\n$t.Cell(1, 5).Range.Text = \"66\u00f78=8, 2\" # '16\u00f73=5, 1' -> '66\u00f78=8, 2'\n$t.Cell(5, 1).Range.Text = \"65\u00f78=8, 1\" # '43\u00f72=21, 1' -> '65\u00f78=8, 1'\n$t.Cell(5, 2).Range.Text = \"40\u00f78=5, 0\" # '59\u00f74=14, 3' -> '40\u00f78=5, 0'\n$t.Cell(5, 3).Range.Text = \"19\u00f79=2, 1\" # '32\u00f75=6, 2' -> '19\u00f79=2, 1'\n$t.Cell(5, 4).Range.Text = \"22\u00f76=3, 4\" # '55\u00f73=18, 1' -> '22\u00f76=3, 4'\n$t.Cell(5, 5).Range.Text = \"63\u00f76=10, 3\" # '53\u00f74=13, 1' -> '63\u00f76=10, 3'\n$t.Cell(9, 1).Range.Text = \"18\u00f79=2, 0\" # '68\u00f74=17, 0' -> '18\u00f79=2, 0'\n$t.Cell(9, 2).Range.Text = \"29\u00f79=3, 2\" # '86\u00f72=43, 0' -> '29\u00f79=3, 2'\n$t.Cell(9, 3).Range.Text = \"32\u00f78=4, 0\" # '85\u00f76=14, 1' -> '32\u00f78=4, 0'\n$t.Cell(9, 4).Range.Text = \"68\u00f73=22, 2\" # '99\u00f78=12, 3' -> '68\u00f73=22, 2'\n$t.Cell(9, 5).Range.Text = \"83\u00f72=41, 1\" # '53\u00f78=6, 5' -> '83\u00f72=41, 1'\n$t.Cell(13, 1).Range.Text = \"99\u00f78=12, 3\" # '60\u00f73=20, 0' -> '99\u00f78=12, 3'\n$t.Cell(13, 2).Range.Text = \"98\u00f79=10, 8\" # '25\u00f76=4, 1' -> '98\u00f79=10, 8'\n$t.Cell(13, 3).Range.Text = \"11\u00f76=1, 5\" # '10\u00f75=2, 0' -> '11\u00f76=1, 5'\n$t.Cell(13, 4).Range.Text = \"44\u00f74=11, 0\" # '57\u00f75=11, 2' -> '44\u00f74=11, 0'\n$t.Cell(13, 5).Range.Text = \"15\u00f76=2, 3\" # '85\u00f76=14, 1' -> '15\u00f76=2, 3'\n$t.Cell(17, 1).Range.Text = \"27\u00f79=3, 0\" # '20\u00f78=2, 4' -> '27\u00f79=3, 0'\n$t.Cell(17, 2).Range.Text = \"53\u00f73=17, 2\" # '13\u00f76=2, 1' -> '53\u00f73=17, 2'\n$t.Cell(17, 3).Range.Text = \"50\u00f72=25, 0\" # '47\u00f75=9, 2' -> '50\u00f72=25, 0'\n$t.Cell(17, 4).Range.Text = \"81\u00f73=27, 0\" # '24\u00f78=3, 0' -> '81\u00f73=27, 0'\n$t.Cell(17, 5).Range.Text = \"71\u00f78=8, 7\" # '28\u00f72=14, 0' -> '71\u00f78=8, 7'\n"}
